# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.007.22"
$ws.Range("E2").Value = "  +1.84%  "
$ws.Range("D3").Value = "3.345.75"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'583.54"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").Value = "'177.95"
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").Value = "3.345.89"
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("E10").Value = "  +7.75%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D12").Value = "'47.54"
$ws.Range("E12").Value = "  +5.92%  "
$ws.Range("E13").Value = "  +3.55%  "
$ws.Range("D14").Value = "'712.51"
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("D15").Value = "3.882.34"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("E16").Value = "  +2.83%  "
$ws.Range("D17").Value = "68.065.02"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "3.351.69"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").Value = "'17.58"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "'11.14"
$ws.Range("E21").Value = "  +4.62%  "
$ws.Range("D22").Value = "'0.900"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").Value = "'5.39"
$ws.Range("E23").Value = "  +3.68%  "
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").Value = "'100.51"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").Value = "'3.93"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").Value = "'2.72"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").Value = "'9.66"
$ws.Range("E28").Value = "  +5.94%  "
$ws.Range("D29").Value = "'33.30"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "'8.63"
$ws.Range("E30").Value = "  +3.91%  "
$ws.Range("D31").Value = "'7.09"
$ws.Range("E31").Value = "  +7.09%  "
$ws.Range("D32").Value = "'570.90"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "'11.05"
$ws.Range("E33").Value = "  +2.66%  "
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'57.77"
$ws.Range("E35").Value = "  +4.94%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").Value = "'3.46"
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").Value = "3.708.79"
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D39").Value = "'34.72"
$ws.Range("E39").Value = "  +10.81%  "
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("D41").Value = "'2.67"
$ws.Range("E41").Value = "  +3.95%  "
$ws.Range("D42").Value = "'3.19"
$ws.Range("E42").Value = "  +7.20%  "
$ws.Range("E43").Value = "  +3.11%  "
$ws.Range("D44").Value = "'0.339"
$ws.Range("E44").Value = "  +4.81%  "
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "'0.0409"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("E47").Value = "  +7.29%  "
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "'131.19"
$ws.Range("E51").Value = "  +1.30%  "
